$wb = $excel.ActiveWorkbook

# Rename sheets (by index to be safe regardless of current name)
$wb.Worksheets.Item(1).Name = "GNG_TO-16502911437090864"
$wb.Worksheets.Item(2).Name = "NB_TO-16502911460622582"
$wb.Worksheets.Item(3).Name = "RS_TO-16502911460632634"
$wb.Worksheets.Item(4).Name = "TOL_TO-16502911461246421"
$wb.Worksheets.Item(5).Name = "vSAT_TO-1650291146217783"

# Sheet 1 (GNG_TO) - update B2:B5
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16502911436559818.csv"
$ws1.Range("B3").Value = "GNG_stims-16502911436772716.csv"
$ws1.Range("B4").Value = "go_stims-1650291143678296.csv"
$ws1.Range("B5").Value = "GNG_stims-16502911437090864.csv"

# Sheet 2 (NB_TO) - update B2:B10
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "ZB-match_5-1650291143896254.csv"
$ws2.Range("B3").Value = "TB-16502911460400543.csv"
$ws2.Range("B4").Value = "OB-16502911447956412.csv"
$ws2.Range("B5").Value = "TB-16502911449843614.csv"
$ws2.Range("B6").Value = "ZB-match_2-16502911439286165.csv"
$ws2.Range("B7").Value = "TB-16502911455180576.csv"
$ws2.Range("B8").Value = "ZB-match_0-165029114384006.csv"
$ws2.Range("B9").Value = "OB-1650291144914571.csv"
$ws2.Range("B10").Value = "OB-16502911446438317.csv"

# Sheet 4 (TOL_TO) - update B2:B7
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-165029114607745.csv"
$ws4.Range("B3").Value = "ZM_stims-16502911460654263.csv"
$ws4.Range("B4").Value = "MM_stims-16502911461089785.csv"
$ws4.Range("B5").Value = "ZM_stims-16502911460784261.csv"
$ws4.Range("B6").Value = "MM_stims-16502911461236389.csv"
$ws4.Range("B7").Value = "ZM_stims-16502911461099434.csv"

# Sheet 5 (vSAT_TO) - update B2:B5
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-1650291146155416.csv"
$ws5.Range("B3").Value = "vSAT_stims-16502911461714208.csv"
$ws5.Range("B4").Value = "vSAT_stims-16502911462017767.csv"
$ws5.Range("B5").Value = "SAT_stims-16502911461292996.csv"
